# Edit slide 9 (the last slide): expand the "GIT HUB :" title into a
# full line of text that includes the GitHub repo link, split across
# several runs with hyperlinks, matching the author's pasted-link edit.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$shp = $s.Shapes.Item(2)          # "제목 4" title placeholder shape
$tr = $shp.TextFrame.TextRange

# Replace the whole paragraph text first so run-splitting via
# Characters() below operates on the final character stream.
$tr.Text = "GIT HUB : YeJin-Choi88/opensw (github.com))"

$repoUrl = "https://github.com/YeJin-Choi88/opensw"
# The trailing ")" was linked separately (a stray paren left over from
# pasting a markdown-style link), ending up pointing at a slightly
# different (malformed) address and thus its own relationship id.
$strayUrl = "https://github.com/YeJin-Choi88/opensw)"

# "YeJin-Choi88/"  -> characters 11-23
$tr.Characters(11, 13).ActionSettings.Item(1).Hyperlink.Address = $repoUrl

# "opensw"         -> characters 24-29
$tr.Characters(24, 6).ActionSettings.Item(1).Hyperlink.Address = $repoUrl

# " (github.com)"  -> characters 30-42
$tr.Characters(30, 13).ActionSettings.Item(1).Hyperlink.Address = $repoUrl

# ")"              -> character 43
$tr.Characters(43, 1).ActionSettings.Item(1).Hyperlink.Address = $strayUrl
